# daily auto push: 2026-01-16 13:47 UTC
#
# A new timestamped data point for 2026/01/16 (金, 19:00, ranking 34) needs to
# be inserted into the log table right after the existing 2026/01/16 entries
# (which currently end at row 636), pushing the 2026/12/29 .. 2027/01/05 block
# down by one row (rows 637:678 -> 638:679).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert one new blank row at 637, shifting rows 637:678 down to 638:679.
$ws.Rows("637:637").Insert(-4121)   # -4121 == xlShiftDown

# The "date" column holds plain text (e.g. "2026/12/29"), not real Excel
# dates, so force the cell to Text format before typing the value -
# otherwise Excel would auto-convert the slash-separated string into a
# date serial number.
$ws.Range("A637").NumberFormat = "@"
$ws.Range("A637").Value = "2026/01/16"
$ws.Range("B637").Value = "金"
$ws.Range("C637").Value = 19
$ws.Range("D637").Value = 34
